# Root cause analysis sheet: insert a "Supplier" column and populate the
# matching-phrase root-cause detail column, per the commit
# "root cause with matching phrases from a different excel sheet".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Detailed root-cause text, written while the sheet still has its
# original column layout (E = Root Cause) -----------------------------------
# The long narrative root cause for the one fully-detailed record,
# left aligned + vertically centered + wrapped, row grown to fit it.
$ws.Range("E2").Value = "Event:`nProcess control`nSupplier didn't add the issue number to the board because they didn't know that they had to do this`n"
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").VerticalAlignment = -4108
$ws.Range("E2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 80.25

# --- 2. Insert a new column B ("Supplier") -------------------------------
# This shifts the old B:F (Material, PO, Summary, Root Cause,
# Corrective and Preventive actions) one place right, to C:G.
$ws.Columns.Item(2).Insert()

# --- 3. Header row ---------------------------------------------------------
$ws.Range("B1").Value = "Supplier"

# --- 4. Supplier letters for each existing data row -------------------------
$ws.Range("B2").Value  = "A"
$ws.Range("B3").Value  = "B"
$ws.Range("B4").Value  = "C"
$ws.Range("B5").Value  = "C"
$ws.Range("B6").Value  = "C"
$ws.Range("B7").Value  = "D"
$ws.Range("B8").Value  = "A"
$ws.Range("B9").Value  = "B"
$ws.Range("B10").Value = "B"
$ws.Range("B11").Value = "A"
$ws.Range("B12").Value = "D"

# --- 5. Remaining detailed root-cause text (column F, was E before insert) -
# Rows 3, 4 and 8: short multi-line "matching phrase" notes, wrapped only.
$ws.Range("F3").Value = "Machine`nHello`nE-data"
$ws.Range("F3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 45

$ws.Range("F4").Value = "Material`nRepair"
$ws.Range("F4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 30

$ws.Range("F8").Value = "Human Error`nRepair"
$ws.Range("F8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 30

# Remaining rows reuse existing root-cause phrases (shared strings).
$ws.Range("F5").Value = "Material"
$ws.Range("F6").Value = "Material"
$ws.Range("F7").Value = "Human Error"
$ws.Range("F9").Value = "Machine"
$ws.Range("F10").Value = "Machine"
$ws.Range("F11").Value = "Machine"
$ws.Range("F12").Value = "Machine"

# --- 5. Column widths --------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.66666666666667   # Supplier
$ws.Columns.Item(6).ColumnWidth = 51.66666666666667   # Root Cause (wider now)

# --- 6. Selection / scrolled view ------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("F9").Select()
